$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the two newly appended rows (48-49) by
# copying the format of an existing data row, then overwrite values below.
$ws.Range("A2:D2").Copy()
$ws.Range("A48:D48").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A49:D49").PasteSpecial(-4122)

# Rewrite rows 2-49 with the final sorted (alphabetical by ANr) data,
# plus the two newly added entries (A_GERMANY, A_EU) at the end.
$ws.Cells.Item(2, 1).Value = 'A_AGE_00-02'
$ws.Cells.Item(2, 2).Value = 'K_AGE'
$ws.Cells.Item(2, 3).Value = '2 Jahre und jünger'
$ws.Cells.Item(2, 4).Value = '2 years and younger'
$ws.Cells.Item(3, 1).Value = 'A_AGE_03-05'
$ws.Cells.Item(3, 2).Value = 'K_AGE'
$ws.Cells.Item(3, 3).Value = '3 bis unter 5 Jahre'
$ws.Cells.Item(3, 4).Value = '3 to under 5 years'
$ws.Cells.Item(4, 1).Value = 'A_AGE_03-11'
$ws.Cells.Item(4, 2).Value = 'K_AGE'
$ws.Cells.Item(4, 3).Value = '3 bis unter 11 Jahre'
$ws.Cells.Item(4, 4).Value = '3 to under 11 years'
$ws.Cells.Item(5, 1).Value = 'A_AGE_11-18'
$ws.Cells.Item(5, 2).Value = 'K_AGE'
$ws.Cells.Item(5, 3).Value = '11 bis unter 18 Jahre'
$ws.Cells.Item(5, 4).Value = '11 to under 18 years'
$ws.Cells.Item(6, 1).Value = 'A_AGE_12-18'
$ws.Cells.Item(6, 2).Value = 'K_AGE'
$ws.Cells.Item(6, 3).Value = '12 bis unter 18 Jahre'
$ws.Cells.Item(6, 4).Value = '12 to under 18 years'
$ws.Cells.Item(7, 1).Value = 'A_AGE_15-999'
$ws.Cells.Item(7, 2).Value = 'K_AGE'
$ws.Cells.Item(7, 3).Value = '15 Jahre und älter'
$ws.Cells.Item(7, 4).Value = '15 years and older'
$ws.Cells.Item(8, 1).Value = 'A_AGE_20-65'
$ws.Cells.Item(8, 2).Value = 'K_AGE'
$ws.Cells.Item(8, 3).Value = '20 bis unter 65 Jahre'
$ws.Cells.Item(8, 4).Value = '20 to under 65 years'
$ws.Cells.Item(9, 1).Value = 'A_AGE_60-65'
$ws.Cells.Item(9, 2).Value = 'K_AGE'
$ws.Cells.Item(9, 3).Value = '60 bis unter 65 Jahre'
$ws.Cells.Item(9, 4).Value = '60 to under 65 years'
$ws.Cells.Item(10, 1).Value = 'A_AIRPOLL_NH3'
$ws.Cells.Item(10, 2).Value = 'K_AIRPOLL'
$ws.Cells.Item(10, 3).Value = 'NH₃'
$ws.Cells.Item(10, 4).Value = 'NH₃'
$ws.Cells.Item(11, 1).Value = 'A_AIRPOLL_NMVOC'
$ws.Cells.Item(11, 2).Value = 'K_AIRPOLL'
$ws.Cells.Item(11, 3).Value = 'NMVOC'
$ws.Cells.Item(11, 4).Value = 'NMVOC'
$ws.Cells.Item(12, 1).Value = 'A_AIRPOLL_NOX'
$ws.Cells.Item(12, 2).Value = 'K_AIRPOLL'
$ws.Cells.Item(12, 3).Value = 'NOₓ'
$ws.Cells.Item(12, 4).Value = 'NOₓ'
$ws.Cells.Item(13, 1).Value = 'A_AIRPOLL_SO2'
$ws.Cells.Item(13, 2).Value = 'K_AIRPOLL'
$ws.Cells.Item(13, 3).Value = 'SO₂'
$ws.Cells.Item(13, 4).Value = 'SO₂'
$ws.Cells.Item(14, 1).Value = 'A_AREA_EU'
$ws.Cells.Item(14, 2).Value = 'K_AREA'
$ws.Cells.Item(14, 3).Value = 'EU'
$ws.Cells.Item(14, 4).Value = 'EU'
$ws.Cells.Item(15, 1).Value = 'A_AREA_EU27'
$ws.Cells.Item(15, 2).Value = 'K_AREA'
$ws.Cells.Item(15, 3).Value = 'EU-27'
$ws.Cells.Item(15, 4).Value = 'EU-27'
$ws.Cells.Item(16, 1).Value = 'A_AREA_EU28'
$ws.Cells.Item(16, 2).Value = 'K_AREA'
$ws.Cells.Item(16, 3).Value = 'EU-28'
$ws.Cells.Item(16, 4).Value = 'EU-28'
$ws.Cells.Item(17, 1).Value = 'A_AREA_EURZ'
$ws.Cells.Item(17, 2).Value = 'K_AREA'
$ws.Cells.Item(17, 3).Value = 'Eurozone'
$ws.Cells.Item(17, 4).Value = 'Eurozone'
$ws.Cells.Item(18, 1).Value = 'A_AREA_GMN'
$ws.Cells.Item(18, 2).Value = 'K_AREA'
$ws.Cells.Item(18, 3).Value = 'Deutschland'
$ws.Cells.Item(18, 4).Value = 'Germany'
$ws.Cells.Item(19, 1).Value = 'A_CRIM_BURGLERY'
$ws.Cells.Item(19, 2).Value = 'K_CRIM'
$ws.Cells.Item(19, 3).Value = 'Wohnungseinbruchsdiebstahl'
$ws.Cells.Item(19, 4).Value = 'Domestic burglary'
$ws.Cells.Item(20, 1).Value = 'A_CRIM_FRAUD'
$ws.Cells.Item(20, 2).Value = 'K_CRIM'
$ws.Cells.Item(20, 3).Value = 'Betrug'
$ws.Cells.Item(20, 4).Value = 'Fraud'
$ws.Cells.Item(21, 1).Value = 'A_CRIM_HARM'
$ws.Cells.Item(21, 2).Value = 'K_CRIM'
$ws.Cells.Item(21, 3).Value = 'Gefährliche und schwere Körperverletzung'
$ws.Cells.Item(21, 4).Value = 'Dangerous and serious bodily injury'
$ws.Cells.Item(22, 1).Value = 'A_CRIM_OTHER'
$ws.Cells.Item(22, 2).Value = 'K_CRIM'
$ws.Cells.Item(22, 3).Value = 'Sonstige Straftaten'
$ws.Cells.Item(22, 4).Value = 'Other offences'
$ws.Cells.Item(23, 1).Value = 'A_LAENDER_BB'
$ws.Cells.Item(23, 2).Value = 'K_LAENDER'
$ws.Cells.Item(23, 3).Value = 'Brandenburg'
$ws.Cells.Item(23, 4).Value = 'Brandenburg'
$ws.Cells.Item(24, 1).Value = 'A_LAENDER_BE'
$ws.Cells.Item(24, 2).Value = 'K_LAENDER'
$ws.Cells.Item(24, 3).Value = 'Berlin'
$ws.Cells.Item(24, 4).Value = 'Berlin'
$ws.Cells.Item(25, 1).Value = 'A_LAENDER_BW'
$ws.Cells.Item(25, 2).Value = 'K_LAENDER'
$ws.Cells.Item(25, 3).Value = 'Baden-Württemberg'
$ws.Cells.Item(25, 4).Value = 'Baden-Wuerttemberg'
$ws.Cells.Item(26, 1).Value = 'A_LAENDER_BY'
$ws.Cells.Item(26, 2).Value = 'K_LAENDER'
$ws.Cells.Item(26, 3).Value = 'Bayern'
$ws.Cells.Item(26, 4).Value = 'Bavaria'
$ws.Cells.Item(27, 1).Value = 'A_LAENDER_HB'
$ws.Cells.Item(27, 2).Value = 'K_LAENDER'
$ws.Cells.Item(27, 3).Value = 'Bremen'
$ws.Cells.Item(27, 4).Value = 'Bremen'
$ws.Cells.Item(28, 1).Value = 'A_LAENDER_HE'
$ws.Cells.Item(28, 2).Value = 'K_LAENDER'
$ws.Cells.Item(28, 3).Value = 'Hessen'
$ws.Cells.Item(28, 4).Value = 'Hesse'
$ws.Cells.Item(29, 1).Value = 'A_LAENDER_HH'
$ws.Cells.Item(29, 2).Value = 'K_LAENDER'
$ws.Cells.Item(29, 3).Value = 'Hamburg'
$ws.Cells.Item(29, 4).Value = 'Hamburg'
$ws.Cells.Item(30, 1).Value = 'A_LAENDER_MV'
$ws.Cells.Item(30, 2).Value = 'K_LAENDER'
$ws.Cells.Item(30, 3).Value = 'Mecklenburg-Vorpommern'
$ws.Cells.Item(30, 4).Value = 'Mecklenburg Western Pomerania'
$ws.Cells.Item(31, 1).Value = 'A_LAENDER_NI'
$ws.Cells.Item(31, 2).Value = 'K_LAENDER'
$ws.Cells.Item(31, 3).Value = 'Niedersachsen'
$ws.Cells.Item(31, 4).Value = 'Lower Saxony'
$ws.Cells.Item(32, 1).Value = 'A_LAENDER_NW'
$ws.Cells.Item(32, 2).Value = 'K_LAENDER'
$ws.Cells.Item(32, 3).Value = 'Nordrhein-Westfalen'
$ws.Cells.Item(32, 4).Value = 'North Rhine-Westphalia'
$ws.Cells.Item(33, 1).Value = 'A_LAENDER_RP'
$ws.Cells.Item(33, 2).Value = 'K_LAENDER'
$ws.Cells.Item(33, 3).Value = 'Rheinland-Pfalz'
$ws.Cells.Item(33, 4).Value = 'Rhineland Palatinate'
$ws.Cells.Item(34, 1).Value = 'A_LAENDER_SH'
$ws.Cells.Item(34, 2).Value = 'K_LAENDER'
$ws.Cells.Item(34, 3).Value = 'Schleswig-Holstein'
$ws.Cells.Item(34, 4).Value = 'Schleswig-Holstein'
$ws.Cells.Item(35, 1).Value = 'A_LAENDER_SL'
$ws.Cells.Item(35, 2).Value = 'K_LAENDER'
$ws.Cells.Item(35, 3).Value = 'Saarland'
$ws.Cells.Item(35, 4).Value = 'Saarland'
$ws.Cells.Item(36, 1).Value = 'A_LAENDER_SN'
$ws.Cells.Item(36, 2).Value = 'K_LAENDER'
$ws.Cells.Item(36, 3).Value = 'Sachsen'
$ws.Cells.Item(36, 4).Value = 'Saxony'
$ws.Cells.Item(37, 1).Value = 'A_LAENDER_ST'
$ws.Cells.Item(37, 2).Value = 'K_LAENDER'
$ws.Cells.Item(37, 3).Value = 'Sachsen-Anhalt'
$ws.Cells.Item(37, 4).Value = 'Saxony-Anhalt'
$ws.Cells.Item(38, 1).Value = 'A_LAENDER_TH'
$ws.Cells.Item(38, 2).Value = 'K_LAENDER'
$ws.Cells.Item(38, 3).Value = 'Thüringen'
$ws.Cells.Item(38, 4).Value = 'Thuringia'
$ws.Cells.Item(39, 1).Value = 'A_PM2.5'
$ws.Cells.Item(39, 2).Value = 'K_PM'
$ws.Cells.Item(39, 3).Value = 'PM2.5'
$ws.Cells.Item(39, 4).Value = 'PM2.5'
$ws.Cells.Item(40, 1).Value = 'A_SEA_B'
$ws.Cells.Item(40, 2).Value = 'K_SEA'
$ws.Cells.Item(40, 3).Value = 'Ostsee'
$ws.Cells.Item(40, 4).Value = 'Baltic sea'
$ws.Cells.Item(41, 1).Value = 'A_SEA_N'
$ws.Cells.Item(41, 2).Value = 'K_SEA'
$ws.Cells.Item(41, 3).Value = 'Nordsee'
$ws.Cells.Item(41, 4).Value = 'Greater North Sea'
$ws.Cells.Item(42, 1).Value = 'A_SEX_D'
$ws.Cells.Item(42, 2).Value = 'K_SEX'
$ws.Cells.Item(42, 3).Value = 'Divers'
$ws.Cells.Item(42, 4).Value = 'Divers'
$ws.Cells.Item(43, 1).Value = 'A_SEX_F'
$ws.Cells.Item(43, 2).Value = 'K_SEX'
$ws.Cells.Item(43, 3).Value = 'Weiblich'
$ws.Cells.Item(43, 4).Value = 'Female'
$ws.Cells.Item(44, 1).Value = 'A_SEX_M'
$ws.Cells.Item(44, 2).Value = 'K_SEX'
$ws.Cells.Item(44, 3).Value = 'Männlich'
$ws.Cells.Item(44, 4).Value = 'Male'
$ws.Cells.Item(45, 1).Value = 'A_SEX_U'
$ws.Cells.Item(45, 2).Value = 'K_SEX'
$ws.Cells.Item(45, 3).Value = 'Unbekannt'
$ws.Cells.Item(45, 4).Value = 'Unknown'
$ws.Cells.Item(46, 1).Value = 'A_URBAN_NONRURAL'
$ws.Cells.Item(46, 2).Value = 'K_URBAN'
$ws.Cells.Item(46, 3).Value = 'Nicht-ländliche Gebiete'
$ws.Cells.Item(46, 4).Value = 'Non-rural areas'
$ws.Cells.Item(47, 1).Value = 'A_URBAN_RURAL'
$ws.Cells.Item(47, 2).Value = 'K_URBAN'
$ws.Cells.Item(47, 3).Value = 'Ländliche Gebiete'
$ws.Cells.Item(47, 4).Value = 'Rural areas'
$ws.Cells.Item(48, 1).Value = 'A_GERMANY'
$ws.Cells.Item(48, 2).Value = 'K_AREA'
$ws.Cells.Item(48, 3).Value = 'Deutschland'
$ws.Cells.Item(48, 4).Value = 'Germany'
$ws.Cells.Item(49, 1).Value = 'A_EU'
$ws.Cells.Item(49, 2).Value = 'K_AREA'
$ws.Cells.Item(49, 3).Value = 'Eropäische Union'
$ws.Cells.Item(49, 4).Value = 'European Union'
